# Fix expense ("Utgift") values: change 2160000 -> 180000 in column E,
# rows 7-72, on both the "private" and "Income" worksheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("private", "Income")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    for ($row = 7; $row -le 72; $row++) {
        $cell = $ws.Cells.Item($row, 5)  # Column E
        if ($cell.Value2 -eq 2160000) {
            $cell.Value = 180000
        }
    }
}
